# Append new translation rows to the "Import" sheet and move the
# selection/scroll position to match the author's final cursor location.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New rows 601..607 -- column A is always "cs", column B is the
# translation key, column C is the Czech translation text.
$rows = @(
    @("cs", "root.common.utils.title", "Nástroje"),
    @("cs", "root.cache.drop.label", "Smazat cache"),
    @("cs", "root.cache.drop.success", "Cache byla smazána"),
    @("cs", "lab.build.button.deactivate", "Deaktivovat build"),
    @("cs", "lab.build.button.activate", "Aktivovat build"),
    @("cs", "lab.build.deactivated.success", "Build [{{data.name}}] byl deaktivován; přestane se nabízet v různých nabídkách napříč aplikací."),
    @("cs", "lab.build.activated.success", "Build [{{data.name}}] byl aktivován; bude se opět nabízet v nabídkách napříč aplikací.")
)

$startRow = 601
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowRange = $ws.Range("A" + $r + ":C" + $r)
    # Match the formatting already used by every other data row ("import"
    # cell style: wrapped text, 10pt Calibri).
    $rowRange.WrapText = $true
    $rowRange.Font.Size = 10

    $ws.Range("A" + $r).Value = $rows[$i][0]
    $ws.Range("B" + $r).Value = $rows[$i][1]
    $ws.Range("C" + $r).Value = $rows[$i][2]
}

# Row 606's translation is long enough to wrap onto a second line in
# Excel, which grows the row to ~2 lines tall.
$ws.Rows.Item(606).RowHeight = 26.25

# Leave the selection where the author ended up.
[void]$ws.Range("B601").Select()
